# Update clinical feature percentages with non-missing N prefix and
# swap counts to reflect the "Latest version of the figures".
$d = $word.ActiveDocument

$replacements = @(
    @{ Old = "40 (80%)";  New = "[50] 10 (20%)" },
    @{ Old = "42 (84%)";  New = "[50] 8 (16%)" },
    @{ Old = "33 (66%)";  New = "[50] 17 (34%)" },
    @{ Old = "28 (56%)";  New = "[50] 22 (44%)" },
    @{ Old = "35 (70%)";  New = "[50] 15 (30%)" },
    @{ Old = "48 (96%)";  New = "[50] 2 (4.0%)" },
    @{ Old = "38 (76%)";  New = "[50] 12 (24%)" },
    @{ Old = "50 (100%)"; New = "[50] 0 (0%)" },
    @{ Old = "n (%)";     New = "[N Non-missing] n (%)" }
)

# Find.Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards,
#              MatchSoundsLike, MatchAllWordForms, Forward, Wrap, Format,
#              ReplaceWith, Replace)
foreach ($pair in $replacements) {
    $range = $d.Content
    $range.Find.Execute($pair.Old, $true, $false, $false, $false, $false,
                         $true, 1, $false, $pair.New, 2)
}
